# Apply changes described by the diff:
# 1. Update selection to A8
# 2. Change B7 and B8 values to "19011200077777"
# 3. Fill in B13..B21 with new serial numbers
# Values must be stored as text (inline strings), not numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "19011200077777"
$ws.Range("B8").Value = "19011200077777"

$ws.Range("B13").Value = "19011200030003"
$ws.Range("B14").Value = "190112000777347"
$ws.Range("B15").Value = "190112000777347"
$ws.Range("B16").Value = "190112000777347"
$ws.Range("B17").Value = "19011200076347"
$ws.Range("B18").Value = "19011200076347"
$ws.Range("B19").Value = "19011200076347"
$ws.Range("B20").Value = "19011200076347"
$ws.Range("B21").Value = "19011200076347"

# Update the active selection to A8 (single cell)
$ws.Range("A8").Select()
